$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.508.01"
$ws.Range("E2").Value = "  +2.34%  "

$ws.Range("D3").Value = "3.810.02"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "681.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.92%  "

$ws.Range("D7").Value = "3.809.00"
$ws.Range("E7").Value = "  +1.16%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  +0.87%  "

$ws.Range("E10").Value = "  +1.80%  "

$ws.Range("E11").Value = "  +7.15%  "

$ws.Range("E12").Value = "  +0.86%  "

$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.68%  "

$ws.Range("D15").Value = "4.452.16"
$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").Value = "3.799.38"
$ws.Range("E16").Value = "  +0.87%  "

$ws.Range("D17").Value = "70.538.40"
$ws.Range("E17").Value = "  +2.35%  "

$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("E19").Value = "  +2.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.115"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +19.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "476.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.39%  "

$ws.Range("E23").Value = "  +1.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.72%  "

$ws.Range("E25").Value = "  -1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.27%  "

$ws.Range("E27").Value = "  +3.32%  "

$ws.Range("E28").Value = "  -1.05%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").Value = "3.961.24"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("E31").Value = "  +9.95%  "

$ws.Range("E32").Value = "  +3.13%  "

$ws.Range("E33").Value = "  +4.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.22%  "

$ws.Range("E35").Value = "  +5.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").Value = "3.759.06"
$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("E39").Value = "  +1.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.44%  "

$ws.Range("E41").Value = "  +2.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.964"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "

$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("E44").Value = "  +12.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000296"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.29%  "

$ws.Range("E51").Value = "  +1.86%  "
